# First addition of 033 R Markdown
# Update the "Module" value for weeks 3-5 on the Schedule_date sheet
# from "2: Coding fundamental" to "2: Coding fundamentals".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule_date")

# Update column C (Module) for rows 4-6 (Week 3, 4, 5)
$ws.Range("C4:C6").Value = "2: Coding fundamentals"

# Move the active selection to D7 as recorded in the saved workbook view
$ws.Activate()
$ws.Range("D7").Select()
